# [DISC-5] Test 페이지 개발 3
#
# Duplicate slide 3 (the "제목 및 내용" slide holding the blue
# rectangle+triangle callout group) to create a new slide that is
# inserted immediately after it (pushing the old slide 4 - the picture
# slide - down to slide 5). Then recolor the duplicated group's four
# shapes from blue (4298B4) to green (33A474), matching the new
# "DISC" callout group, and restore the expected Korean group name.

$p = $ppt.ActivePresentation

$srcSlide = $p.Slides.Item(3)

# Duplicate() inserts the copy right after the source slide and shifts
# everything after it down by one - exactly the sldIdLst change we need
# (..., 258, 260[new], 259[old slide4], ...).
$newRange = $srcSlide.Duplicate()
$newSlide = $newRange.Item(1)

$grp = $newSlide.Shapes.Item(1)
$grp.Name = "그룹 1"

$greenRGB = 7644211  # srgbClr 33A474 (COM RGB is 0xBBGGRR)

for ($i = 1; $i -le $grp.GroupItems.Count; $i++) {
    $shp = $grp.GroupItems.Item($i)
    $shp.Fill.ForeColor.RGB = $greenRGB
}
